$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Primers_for_verify2")

# Insert a new row at position 3 (shifts old row3 -> row4, old row4 -> row5)
$ws.Rows.Item(3).Insert()

# Match the header-style formatting used by the other "ID" cells in column A
# (bold font, centered/top alignment, thin box border) so the new A3 cell
# renders the same as A2/A4/A5.
$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").HorizontalAlignment = -4108
$ws.Range("A3").VerticalAlignment = -4160
$ws.Range("A3").Borders.LineStyle = 1

# Populate the newly inserted row 3 with the aceE_del primer record
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "aceE_del"
$ws.Range("C3").Value = "AGAAGATGTTGTAAATCAAGCGCATATAAAAGCGCGGCAACTAAACGTAGAACCTGTCTTATTGAGCTTTCCGGCGAGAGTTCAATGGGACAGGTTCCAGAAAACTCAACGTTATTAGATAGATAAGGAATAACCCGAGGTAAAAGAATAATGGCTATCGAAATCAAAGTACCGGACATCGGGGCTGATGAAGTTGAAATCACCGAGATCCTGGTCAAAGTGGGCGACAAAGTTGAAGCCGAACAGTCGCTGATCACCGTAGAAGGCGAC"
$ws.Range("D3").Value = 9.333549347603878
$ws.Range("E3").Value = -22.75753292142625
$ws.Range("F3").Value = -18.03345208828691
$ws.Range("G3").Value = "AGAAGATGTTGTAAATCAAGCGCA"
$ws.Range("H3").Value = "GTCGCCTTCTACGGTGATCA"
$ws.Range("I3").Value = 59.54446661842991
$ws.Range("J3").Value = 59.5453996763095
$ws.Range("K3").Value = 270

# The record that used to be row 4 (Cgl1452_ins) is now row 5; renumber its
# ID from 2 to 3 to keep the sequence (0,1,2,3) consistent.
$ws.Range("A5").Value = 3
